# update scripts wuth new tpm
# Refresh the NATMI ligand-receptor (Ccl25-Ccr10) metrics for each
# sending/target cluster pair with the newly recomputed TPM-based values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 4.055275000000001
$ws.Range("H2").Value = 12.165825
$ws.Range("I2").Value = 0.1947228515851206
$ws.Range("J2").Value = 0.1947228515851206
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.058783666666667
$ws.Range("N2").Value = 3.176351
$ws.Range("O2").Value = 0.5014862149947701
$ws.Range("P2").Value = 0.5014862149947702
$ws.Range("Q2").Value = 4.293658933841668
$ws.Range("R2").Value = 38.64293040457501
$ws.Range("S2").Value = 0.09765082581441051
$ws.Range("T2").Value = 0.09765082581441051

# Row 3
$ws.Range("G3").Value = 4.055275000000001
$ws.Range("H3").Value = 12.165825
$ws.Range("I3").Value = 0.1947228515851206
$ws.Range("J3").Value = 0.1947228515851206
$ws.Range("O3").Value = 0.475327031240749
$ws.Range("P3").Value = 0.4753270312407492
$ws.Range("Q3").Value = 4.069687447350001
$ws.Range("R3").Value = 36.62718702615
$ws.Range("S3").Value = 0.09255703495868838
$ws.Range("T3").Value = 0.09255703495868838

# Row 4
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("G4").Value = 4.055275000000001
$ws.Range("H4").Value = 12.165825
$ws.Range("I4").Value = 0.1947228515851206
$ws.Range("J4").Value = 0.1947228515851206
$ws.Range("M4").Value = 0.048954
$ws.Range("N4").Value = 0.146862
$ws.Range("O4").Value = 0.02318675376448066
$ws.Range("P4").Value = 0.02318675376448067
$ws.Range("Q4").Value = 0.19852193235
$ws.Range("R4").Value = 1.78669739115
$ws.Range("S4").Value = 0.004514990812021705
$ws.Range("T4").Value = 0.004514990812021705

# Row 5
$ws.Range("G5").Value = 6.542812333333333
$ws.Range("I5").Value = 0.3141673684110111
$ws.Range("J5").Value = 0.3141673684110111
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.058783666666667
$ws.Range("N5").Value = 3.176351
$ws.Range("O5").Value = 0.5014862149947701
$ws.Range("P5").Value = 0.5014862149947702
$ws.Range("Q5").Value = 6.927422832598556
$ws.Range("R5").Value = 62.346805493387
$ws.Range("S5").Value = 0.1575506044593055
$ws.Range("T5").Value = 0.1575506044593055

# Row 6
$ws.Range("G6").Value = 6.542812333333333
$ws.Range("I6").Value = 0.3141673684110111
$ws.Range("J6").Value = 0.3141673684110111
$ws.Range("O6").Value = 0.475327031240749
$ws.Range("P6").Value = 0.4753270312407492
$ws.Range("Q6").Value = 6.566065488366
$ws.Range("R6").Value = 59.09458939529399
$ws.Range("S6").Value = 0.1493322425395246
$ws.Range("T6").Value = 0.1493322425395246

# Row 7
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("G7").Value = 6.542812333333333
$ws.Range("I7").Value = 0.3141673684110111
$ws.Range("J7").Value = 0.3141673684110111
$ws.Range("M7").Value = 0.048954
$ws.Range("N7").Value = 0.146862
$ws.Range("O7").Value = 0.02318675376448066
$ws.Range("P7").Value = 0.02318675376448067
$ws.Range("Q7").Value = 0.320296834966
$ws.Range("R7").Value = 2.882671514694
$ws.Range("S7").Value = 0.007284521412180995
$ws.Range("T7").Value = 0.007284521412180996

# Row 8
$ws.Range("G8").Value = 6.7998
$ws.Range("H8").Value = 20.3994
$ws.Range("I8").Value = 0.3265071903159472
$ws.Range("J8").Value = 0.3265071903159472
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.058783666666667
$ws.Range("N8").Value = 3.176351
$ws.Range("O8").Value = 0.5014862149947701
$ws.Range("P8").Value = 0.5014862149947702
$ws.Range("Q8").Value = 7.199517176600001
$ws.Range("R8").Value = 64.7956545894
$ws.Range("S8").Value = 0.1637388550401214
$ws.Range("T8").Value = 0.1637388550401215

# Row 9
$ws.Range("G9").Value = 6.7998
$ws.Range("H9").Value = 20.3994
$ws.Range("I9").Value = 0.3265071903159472
$ws.Range("J9").Value = 0.3265071903159472
$ws.Range("O9").Value = 0.475327031240749
$ws.Range("P9").Value = 0.4753270312407492
$ws.Range("Q9").Value = 6.823966489200001
$ws.Range("R9").Value = 61.4156984028
$ws.Range("S9").Value = 0.1551976934516374
$ws.Range("T9").Value = 0.1551976934516375

# Row 10
$ws.Range("D10").Value = "Resolving-Mac"
$ws.Range("G10").Value = 6.7998
$ws.Range("H10").Value = 20.3994
$ws.Range("I10").Value = 0.3265071903159472
$ws.Range("J10").Value = 0.3265071903159472
$ws.Range("M10").Value = 0.048954
$ws.Range("N10").Value = 0.146862
$ws.Range("O10").Value = 0.02318675376448066
$ws.Range("P10").Value = 0.02318675376448067
$ws.Range("Q10").Value = 0.3328774092
$ws.Range("R10").Value = 2.9958966828
$ws.Range("S10").Value = 0.007570641824188293
$ws.Range("T10").Value = 0.007570641824188294

# Row 11
$ws.Range("G11").Value = 3.427994
$ws.Range("H11").Value = 10.283982
$ws.Range("I11").Value = 0.164602589687921
$ws.Range("J11").Value = 0.164602589687921
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 1.058783666666667
$ws.Range("N11").Value = 3.176351
$ws.Range("O11").Value = 0.5014862149947701
$ws.Range("P11").Value = 0.5014862149947702
$ws.Range("Q11").Value = 3.629504056631334
$ws.Range("R11").Value = 32.66553650968201
$ws.Range("S11").Value = 0.08254592968093268
$ws.Range("T11").Value = 0.08254592968093269

# Row 12
$ws.Range("G12").Value = 3.427994
$ws.Range("H12").Value = 10.283982
$ws.Range("I12").Value = 0.164602589687921
$ws.Range("J12").Value = 0.164602589687921
$ws.Range("O12").Value = 0.475327031240749
$ws.Range("P12").Value = 0.4753270312407492
$ws.Range("Q12").Value = 3.440177090676
$ws.Range("R12").Value = 30.961593816084
$ws.Range("S12").Value = 0.07824006029089861
$ws.Range("T12").Value = 0.07824006029089864

# Row 13
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("G13").Value = 3.427994
$ws.Range("H13").Value = 10.283982
$ws.Range("I13").Value = 0.164602589687921
$ws.Range("J13").Value = 0.164602589687921
$ws.Range("M13").Value = 0.048954
$ws.Range("N13").Value = 0.146862
$ws.Range("O13").Value = 0.02318675376448066
$ws.Range("P13").Value = 0.02318675376448067
$ws.Range("Q13").Value = 0.167814018276
$ws.Range("R13").Value = 1.510326164484
$ws.Range("S13").Value = 0.003816599716089668
$ws.Range("T13").Value = 0.003816599716089669
